$wb = $excel.ActiveWorkbook

# Rename sheets (task order identifiers updated to newer timestamps)
$wb.Worksheets.Item(1).Name = "GNG_TO-165029112392062"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911260204015"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911260214045"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650291126078502"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502911261397064"

# Sheet 1: GNG_TO
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911238633564.csv"
$ws1.Range("B3").Value = "GNG_stims-165029112388897.csv"
$ws1.Range("B4").Value = "go_stims-1650291123896218.csv"
$ws1.Range("B5").Value = "GNG_stims-1650291123919616.csv"

# Sheet 2: NB_TO
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-1650291124794414.csv"
$ws2.Range("B3").Value = "TB-16502911260086472.csv"
$ws2.Range("B4").Value = "TB-1650291125280753.csv"
$ws2.Range("B5").Value = "OB-1650291124661566.csv"
$ws2.Range("B6").Value = "ZB-match_2-16502911244722397.csv"
$ws2.Range("B7").Value = "ZB-match_0-16502911241647518.csv"
$ws2.Range("B8").Value = "OB-16502911247676883.csv"
$ws2.Range("B9").Value = "TB-16502911255252764.csv"
$ws2.Range("B10").Value = "ZB-match_4-16502911243968425.csv"

# Sheet 3: RS_TO
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4: TOL_TO
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650291126036338.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911260234077.csv"
$ws4.Range("B4").Value = "MM_stims-1650291126062146.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911260373404.csv"
$ws4.Range("B6").Value = "MM_stims-16502911260775266.csv"
$ws4.Range("B7").Value = "ZM_stims-165029112606315.csv"

# Sheet 5: vSAT_TO
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650291126092713.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502911261082945.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502911261245873.csv"
$ws5.Range("B5").Value = "SAT_stims-1650291126080496.csv"
